$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "path"
$ws.Range("C1").Value = "alias"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "card."
$ws.Range("F1").Value = "stereotype"
$ws.Range("G1").Value = "id"
$ws.Range("H1").Value = "definition"
$ws.Range("I1").Value = "definitioncode"
